$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 243, shifting existing rows 243:358 down to 244:359.
$ws.Rows(243).Insert()

# Populate the newly inserted row 243 with the new record's data.
$ws.Range("A243").Value = 6
$ws.Range("B243").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C243").Value = "Metropolitana"
$ws.Range("D243").Value = "2023-03-31"
$ws.Range("E243").Value = 13
$ws.Range("F243").Value = 100112026
$ws.Range("G243").Value = "Haba"
$ws.Range("H243").Value = "Sin especificar"
$ws.Range("I243").Value = "Primera"
$ws.Range("J243").Value = 800
$ws.Range("K243").Value = 13000
$ws.Range("L243").Value = 14000
$ws.Range("M243").Value = 13562
$ws.Range("N243").Value = "$/saco 25 kilos"
$ws.Range("O243").Value = "Región Metropolitana"
$ws.Range("P243").Value = 542
$ws.Range("Q243").Value = 25
$ws.Range("R243").Value = "Hortaliza"
